# Update "想去人数" (interested-count) figures in the F column across the
# 展览 / 演出 / 全部类型 sheets to match the freshly scraped snapshot.

$wb = $excel.ActiveWorkbook

# --- 展览 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3
$ws.Range("F6").Value = 930
$ws.Range("F7").Value = 166
$ws.Range("F8").Value = 960
$ws.Range("F9").Value = 749
$ws.Range("F10").Value = 199
$ws.Range("F13").Value = 788
$ws.Range("F14").Value = 258
$ws.Range("F15").Value = 556
$ws.Range("F17").Value = 1304
$ws.Range("F19").Value = 431
$ws.Range("F20").Value = 1118
$ws.Range("F21").Value = 2804
$ws.Range("F22").Value = 1316
$ws.Range("F23").Value = 659
$ws.Range("F24").Value = 168
$ws.Range("F25").Value = 1248
$ws.Range("F27").Value = 975
$ws.Range("F29").Value = 1347
$ws.Range("F30").Value = 33
$ws.Range("F31").Value = 3
$ws.Range("F32").Value = 1341

# --- 演出 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 352
$ws.Range("F9").Value = 37
$ws.Range("F10").Value = 151
$ws.Range("F11").Value = 23

# --- 全部类型 ----------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3
$ws.Range("F10").Value = 352
$ws.Range("F13").Value = 930
$ws.Range("F14").Value = 166
$ws.Range("F16").Value = 960
$ws.Range("F17").Value = 749
$ws.Range("F18").Value = 199
$ws.Range("F21").Value = 37
$ws.Range("F22").Value = 151
$ws.Range("F23").Value = 23
$ws.Range("F26").Value = 788
$ws.Range("F27").Value = 258
$ws.Range("F28").Value = 556
$ws.Range("F30").Value = 1304
$ws.Range("F32").Value = 431
$ws.Range("F33").Value = 1118
$ws.Range("F34").Value = 2804
$ws.Range("F35").Value = 1316
$ws.Range("F36").Value = 659
$ws.Range("F37").Value = 168
$ws.Range("F38").Value = 1248
$ws.Range("F42").Value = 975
$ws.Range("F44").Value = 1347
$ws.Range("F45").Value = 33
$ws.Range("F46").Value = 3
$ws.Range("F47").Value = 1341
